# The authored change swaps the OOXML content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml ("Office Theme") becomes the "Integral"
# theme content that used to live in theme2.xml, and theme2.xml becomes the
# "Office Theme" content that used to live in theme1.xml. The fontScheme and
# fmtScheme blocks are byte-identical between the two themes, so the only
# real content delta is the clrScheme (12 colours + the name attributes).
#
# This COM-interop host only exposes one (document-wide) theme through the
# PowerPoint object model - Master.Theme / NotesMaster.Theme / Design.Theme
# all resolve to the single active theme part (ppt/theme/theme2.xml, the one
# wired up via the slide master / presentation relationship). We drive that
# theme's colour scheme to the values the target diff expects to land in
# that file ("Office Theme"'s palette).

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette: the "Office Theme" clrScheme (previously theme1.xml),
# in the same dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order PowerPoint
# uses for ThemeColorScheme.Item(1..12).
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = HexToRgb $officeColors[$i - 1]
}
